# Update "Generate Report for Handback" timestamps / priority flag.
#
# The underlying edit only rewrites five text values inside the workbook's
# shared-string table (in place, same slots/count):
#   "2016-08-26 08:19:22" -> "2016-08-26 08:20:17"
#   "ht"                  -> "mt"
#   "2016-08-26 08:19:18" -> "2016-08-26 08:19:59"
#   "2016-08-26 08:19:33" -> "2016-08-26 08:20:32"
#   "2016-08-26 08:19:40" -> "2016-08-26 08:20:39"
# Every cell across the three sheets that held one of those exact strings
# must end up with the corresponding new text, including the cells for the
# "916e62c6..." row (row 3) and the "a485f69c..." row (row 4), which
# originally happened to share the very same values.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G), rows 3 & 4
# both originally read "2016-08-26 08:19:22".
$wsOverview.Range("G3").Value = "2016-08-26 08:20:17"
$wsOverview.Range("G4").Value = "2016-08-26 08:20:17"

# zh-cn sheet, rows 3 & 4:
#   Column E = "Priority" (ht -> mt)
#   Column H = "Correspond Handoff Datetime"
#   Column K = "Correspond Handback DateTime"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-26 08:19:59"
$wsZhCn.Range("H4").Value = "2016-08-26 08:19:59"
$wsZhCn.Range("K3").Value = "2016-08-26 08:20:32"
$wsZhCn.Range("K4").Value = "2016-08-26 08:20:32"

# de-de sheet, rows 3 & 4:
#   Column E = "Priority" (ht -> mt) -- shares the same string as zh-cn's E
#   Column H = "Correspond Handoff Datetime" -- shares the same string as
#              Overview's G column
#   Column K = "Correspond Handback DateTime"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-26 08:20:17"
$wsDeDe.Range("H4").Value = "2016-08-26 08:20:17"
$wsDeDe.Range("K3").Value = "2016-08-26 08:20:39"
$wsDeDe.Range("K4").Value = "2016-08-26 08:20:39"
